# "Ajuste na frases de Bandido"
# Update the "Bandido" sheet: row 11's stolen-item / phrase pair changes
# from "a carteira" / "não é mais teu" to "a bike" / "não é mais".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bandido")

# Write D11 first, then C11, so the new shared strings are appended in the
# same order as the source edit ("não é mais" before "a bike").
$ws.Range("D11").Value2 = "não é mais"
$ws.Range("C11").Value2 = "a bike"

# Move the sheet's active-cell selection to C12 (matches the saved
# selection state captured in the workbook).
$ws.Range("C12").Select()

# Reposition the workbook window (cosmetic window-position metadata).
$win = $wb.Windows.Item(1)
$win.Left = 5700
